# The commit introduces a new (out-of-repo) PDF-to-image/OCR experiment;
# the tracked Word fixture itself only picks up the bookkeeping churn that
# real Word leaves behind when it opens/resaves the file: a "_GoBack"
# bookmark is left at the last edit position in the (still empty) body
# paragraph.
#
# Reproduce that with the documented Bookmarks API: wrap a zero-length
# range at the start of the lone paragraph in a "_GoBack" bookmark, which
# serialises to exactly:
#   <w:p>
#     <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#     <w:bookmarkEnd w:id="0"/>
#   </w:p>

$d = $word.ActiveDocument

$targetPara = $d.Paragraphs(1)
$r = $targetPara.Range
$r.Collapse(1)   # wdCollapseStart -- zero-length range, before the pilcrow

# Remove a pre-existing _GoBack bookmark first (idempotent / re-run safe).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Bookmarks.Add("_GoBack", $r)
